$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values ---
# Modelo value changed from "Kwid" to "KWID 1.0"
$ws.Range("G2").Value = "KWID 1.0"

# KM atual do veiculo: 5000 -> 500
$ws.Range("J2").Value = 500

# Data de emplacamento do veiculo: 2018-01-01 -> 2018-03-28
$ws.Range("K2").Value = "3/28/2018"

# --- Add new columns at the end: CNPJ, Razao Social ---
$ws.Range("AC1").Value = "CNPJ"
$ws.Range("AD1").Value = "Razão Social"

# --- Column D width ---
# (target stored width is 20.85546875; the COM width setter here quantizes to
# the nearest 1/6 character, so 20 is the closest achievable input)
$ws.Columns("D").ColumnWidth = 20

# --- View state: scroll to Q1, select Y2 ---
$ws.Range("Y2").Select()
$excel.ActiveWindow.ScrollColumn = 17
